$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "W+,Z"
$ws.Range("G3").Value = "W+,Z"
$ws.Range("G4").Value = "W+,Z"

$ws.Range("H11").Select()
